$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.692.79'
$ws.Range("E2").Value = '  -0.95%  '
$ws.Range("D3").Value = '3.087.85'
$ws.Range("E3").Value = '  -2.36%  '
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.55'
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.32'
$ws.Range("E6").Value = '  +6.51%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("E8").Value = '  +1.99%  '
$ws.Range("D9").Value = '3.086.53'
$ws.Range("E9").Value = '  -1.95%  '
$ws.Range("E10").Value = '  -4.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.85'
$ws.Range("E11").Value = '  -1.04%  '
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.42'
$ws.Range("E13").Value = '  -0.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000240'
$ws.Range("E14").Value = '  -3.52%  '
$ws.Range("D15").Value = '3.600.02'
$ws.Range("E15").Value = '  -2.30%  '
$ws.Range("E16").Value = '  -1.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.17'
$ws.Range("E17").Value = '  -1.39%  '
$ws.Range("D18").Value = '63.678.20'
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("D19").Value = '3.087.18'
$ws.Range("E19").Value = '  -2.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '478.78'
$ws.Range("E20").Value = '  +2.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.58'
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.712'
$ws.Range("E22").Value = '  -3.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.56'
$ws.Range("E23").Value = '  -0.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.43'
$ws.Range("E24").Value = '  +2.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.29'
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.87'
$ws.Range("E26").Value = '  -2.31%  '
$ws.Range("E27").Value = '  +4.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.54'
$ws.Range("E30").Value = '  -1.70%  '
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.17'
$ws.Range("E32").Value = '  -2.54%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("E34").Value = '  -1.88%  '
$ws.Range("D35").Value = '0.0₃0852'
$ws.Range("E35").Value = '  -2.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.47'
$ws.Range("E36").Value = '  +6.30%  '
$ws.Range("E37").Value = '  -1.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.06'
$ws.Range("E38").Value = '  -1.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.21'
$ws.Range("E39").Value = '  -3.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.37'
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.81'
$ws.Range("E41").Value = '  -1.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '446.67'
$ws.Range("E42").Value = '  -4.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.11'
$ws.Range("E43").Value = '  +5.42%  '
$ws.Range("E44").Value = '  -2.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0362'
$ws.Range("E45").Value = '  -3.21%  '
$ws.Range("E46").Value = '  +3.69%  '
$ws.Range("D47").Value = '2.826.29'
$ws.Range("E47").Value = '  -2.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.42'
$ws.Range("E48").Value = '  -1.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.49'
$ws.Range("E49").Value = '  +5.79%  '
$ws.Range("E51").Value = '  +0.99%  '
